$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new header "baseline" in F1, copying header style from E1
$ws.Range("F1").Value = "baseline"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Update existing data rows with new values and add column F
$ws.Range("B2").Value = 0.797
$ws.Range("C2").Value = 0.792
$ws.Range("D2").Value = 0.797
$ws.Range("E2").Value = 0.797
$ws.Range("F2").Value = 0.756

$ws.Range("B3").Value = 0.748
$ws.Range("C3").Value = 0.27
$ws.Range("D3").Value = 0.74
$ws.Range("E3").Value = 0.738
$ws.Range("F3").Value = 0.743

$ws.Range("B4").Value = 0.782
$ws.Range("C4").Value = 0.782
$ws.Range("D4").Value = 0.771
$ws.Range("E4").Value = 0.795
$ws.Range("F4").Value = 0.743

# Add new rows 5 and 6
$ws.Range("A5").Value = "Embeddings"
$ws.Range("B5").Value = 0.801
$ws.Range("C5").Value = 0.63
$ws.Range("D5").Value = 0.778
$ws.Range("E5").Value = 0.786
$ws.Range("F5").Value = 0.742

$ws.Range("A6").Value = "Embeddings + Attr"
$ws.Range("B6").Value = 0.755
$ws.Range("C6").Value = 0.269
$ws.Range("D6").Value = 0.755
$ws.Range("E6").Value = 0.775
$ws.Range("F6").Value = 0.747
